$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.582.48"
$ws.Range("E2").Value = "  +4.02%  "
$ws.Range("D3").Value = "3.004.99"
$ws.Range("E3").Value = "  +4.35%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "507.90"
$ws.Range("E5").Value = "  +7.44%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.19"
$ws.Range("E6").Value = "  +8.38%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  +7.19%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.60"
$ws.Range("E9").Value = "  +14.64%  "
$ws.Range("E10").Value = "  +12.53%  "
$ws.Range("E11").Value = "  +7.09%  "
$ws.Range("E12").Value = "  +4.79%  "
$ws.Range("D13").Value = "3.521.48"
$ws.Range("E13").Value = "  +4.89%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.61"
$ws.Range("E14").Value = "  +10.39%  "
$ws.Range("E15").Value = "  +15.08%  "
$ws.Range("D16").Value = "56.600.40"
$ws.Range("E16").Value = "  +4.11%  "
$ws.Range("D17").Value = "3.006.24"
$ws.Range("E17").Value = "  +4.90%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.82"
$ws.Range("E18").Value = "  +8.04%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.51"
$ws.Range("E19").Value = "  +8.74%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.86"
$ws.Range("E20").Value = "  +11.05%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "327.21"
$ws.Range("E21").Value = "  +9.84%  "
$ws.Range("E22").Value = "  -0.25%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.479"
$ws.Range("E23").Value = "  +8.13%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "62.51"
$ws.Range("E24").Value = "  +6.45%  "
$ws.Range("E25").Value = "  +9.75%  "
$ws.Range("E26").Value = "  +0.50%  "
$ws.Range("D27").Value = "0.0₃0922"
$ws.Range("E27").Value = "  +13.96%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.58"
$ws.Range("E28").Value = "  +6.75%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.03"
$ws.Range("E29").Value = "  +13.55%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.25"
$ws.Range("E30").Value = "  +10.66%  "
$ws.Range("E31").Value = "  +9.29%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.68"
$ws.Range("E32").Value = "  +9.38%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "156.61"
$ws.Range("E33").Value = "  +15.10%  "
$ws.Range("E34").Value = "  +7.50%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.66"
$ws.Range("E35").Value = "  +4.44%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.28"
$ws.Range("E36").Value = "  +4.21%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0676"
$ws.Range("E37").Value = "  +8.76%  "
$ws.Range("E38").Value = "  +2.29%  "
$ws.Range("D39").Value = "3.039.19"
$ws.Range("E39").Value = "  +4.84%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.64"
$ws.Range("E40").Value = "  +5.22%  "
$ws.Range("E41").Value = "  +0.20%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.649"
$ws.Range("E42").Value = "  +8.01%  "
$ws.Range("D43").Value = "2.267.96"
$ws.Range("E43").Value = "  +11.18%  "
$ws.Range("E44").Value = "  +5.04%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.41"
$ws.Range("E45").Value = "  +6.63%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.62"
$ws.Range("E46").Value = "  +6.39%  "
$ws.Range("E47").Value = "  +23.18%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0236"
$ws.Range("E48").Value = "  +9.78%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.80"
$ws.Range("E49").Value = "  +7.57%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.21"
$ws.Range("E50").Value = "  +6.88%  "
$ws.Range("E51").Value = "  +10.84%  "
